$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 22:58"

# Update country rows whose ranking / case numbers changed (COVID data refresh).
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5831043
$ws.Range("C4").Value = 34316
$ws.Range("D4").Value = 3138469
$ws.Range("E4").Value = 2512537
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 837
$ws.Range("H4").Value = 180037

$ws.Range("A8").Value = "Sudafrica"
$ws.Range("B8").Value = 607045
$ws.Range("C8").Value = 3707
$ws.Range("D8").Value = 504127
$ws.Range("E8").Value = 89931
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 144
$ws.Range("H8").Value = 12987

$ws.Range("A23").Value = "Alemania"
$ws.Range("B23").Value = 233850
$ws.Range("C23").Value = 829
$ws.Range("D23").Value = 208950
$ws.Range("E23").Value = 15569
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 9331

$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 101933
$ws.Range("C33").Value = 1217
$ws.Range("D33").Value = 78651
$ws.Range("E33").Value = 22463
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 819

$ws.Range("A76").Value = "Estado de Palestina"
$ws.Range("B76").Value = 18476
$ws.Range("C76").Value = 163
$ws.Range("D76").Value = 11103
$ws.Range("E76").Value = 7248
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 125

$ws.Range("A78").Value = "Costa de Marfil"
$ws.Range("B78").Value = 17374
$ws.Range("C78").Value = 64
$ws.Range("D78").Value = 15106
$ws.Range("E78").Value = 2155
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 113

$ws.Range("A94").Value = "Guayana Francesa"
$ws.Range("B94").Value = 8797
$ws.Range("C94").Value = 20
$ws.Range("D94").Value = 8307
$ws.Range("E94").Value = 435
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 55

$ws.Range("A103").Value = "Mauritania"
$ws.Range("B103").Value = 6894
$ws.Range("C103").Value = 9
$ws.Range("D103").Value = 6203
$ws.Range("E103").Value = 533
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 158

$ws.Range("A107").Value = "Malaui"
$ws.Range("B107").Value = 5382
$ws.Range("C107").Value = 60
$ws.Range("D107").Value = 2998
$ws.Range("E107").Value = 2216
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 2
$ws.Range("H107").Value = 168

$ws.Range("A108").Value = "Republica de Yibuti"
$ws.Range("B108").Value = 5382
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 5233
$ws.Range("E108").Value = 89
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 60

$ws.Range("A115").Value = "Suazilandia"
$ws.Range("B115").Value = 4189
$ws.Range("C115").Value = 61
$ws.Range("D115").Value = 2762
$ws.Range("E115").Value = 1344
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 83

$ws.Range("A123").Value = "Somalia"
$ws.Range("B123").Value = 3269
$ws.Range("C123").Value = 4
$ws.Range("D123").Value = 2396
$ws.Range("E123").Value = 780
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 93

$ws.Range("A126").Value = "Ruanda"
$ws.Range("B126").Value = 2889
$ws.Range("C126").Value = 109
$ws.Range("D126").Value = 1754
$ws.Range("E126").Value = 1124
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 11

$ws.Range("A127").Value = "Tunez"
$ws.Range("B127").Value = 2738
$ws.Range("C127").Value = 131
$ws.Range("D127").Value = 1434
$ws.Range("E127").Value = 1236
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 4
$ws.Range("H127").Value = 68

$ws.Range("A128").Value = "Mali"
$ws.Range("B128").Value = 2699
$ws.Range("C128").Value = 11
$ws.Range("D128").Value = 2010
$ws.Range("E128").Value = 564
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 125

$ws.Range("A129").Value = "Eslovenia"
$ws.Range("B129").Value = 2617
$ws.Range("C129").Value = 43
$ws.Range("D129").Value = 2079
$ws.Range("E129").Value = 407
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 131

$ws.Range("A137").Value = "Angola"
$ws.Range("B137").Value = 2134
$ws.Range("C137").Value = 66
$ws.Range("D137").Value = 814
$ws.Range("E137").Value = 1226
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 94

$ws.Range("A138").Value = "Benin"
$ws.Range("B138").Value = 2095
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 1705
$ws.Range("E138").Value = 351
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 39

$ws.Range("A140").Value = "Sierra Leona"
$ws.Range("B140").Value = 1980
$ws.Range("C140").Value = 8
$ws.Range("D140").Value = 1545
$ws.Range("E140").Value = 366
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 69

$ws.Range("A146").Value = "Aruba"
$ws.Range("B146").Value = 1534
$ws.Range("C146").Value = 70
$ws.Range("D146").Value = 446
$ws.Range("E146").Value = 1081
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 7

$ws.Range("A147").Value = "Uruguay"
$ws.Range("B147").Value = 1516
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 1249
$ws.Range("E147").Value = 225
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 42

$ws.Range("A155").Value = "Togo"
$ws.Range("B155").Value = 1275
$ws.Range("C155").Value = 36
$ws.Range("D155").Value = 903
$ws.Range("E155").Value = 345
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 27
